$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells remain stored as text, matching the
# original inline-string typing, instead of being auto-coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "39.435.32"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "2.157.64"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "227.93"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("D7").Value = "63.84"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.395"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").Value = "0.0854"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "15.98"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "2.478.02"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").Value = "22.09"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "0.812"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "2.164.17"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "39.398.90"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").Value = "71.82"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").Value = "230.27"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.34"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "172.35"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  +2.92%  "
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "2.60"
$ws.Range("E31").Value = "  +4.98%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "4.59"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("D34").Value = "7.13"
$ws.Range("E34").Value = "  +8.11%  "
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "3.55"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "103.46"
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D41").Value = "0.0231"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").Value = "17.59"
$ws.Range("E42").Value = "  -4.09%  "
$ws.Range("D43").Value = "1.529.27"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("E44").Value = "  +3.97%  "
$ws.Range("E45").Value = "  +5.09%  "
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("E48").Value = "  +5.09%  "
$ws.Range("D49").Value = "7.66"
$ws.Range("E49").Value = "  -1.57%  "
$ws.Range("D50").Value = "2.361.45"
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("E51").Value = "  -0.19%  "
